# data provider utils updated
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)   # DATA

# --- RUNMANAGER sheet: flip the "execute" flag for the first two test cases ---
$ws1.Range("C2").Value = "no"
$ws1.Range("C3").Value = "no"

# --- DATA sheet: flip "execute" flag for matching rows, and switch browser on row 8 ---
$ws2.Range("B2").Value = "no"
$ws2.Range("B4").Value = "no"
$ws2.Range("C8").Value = "edge"

# Remove the now-redundant last data row (row 9)
$ws2.Rows.Item(9).Delete()

# --- Update the remembered selections on each sheet ---
$ws1.Range("C3").Select()
$ws2.Range("B8").Select()

# Restore DATA as the active/selected sheet (as it was originally)
$ws2.Select()
